$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "..., що є основою CASE, - засоби." -> "..., що є основою CASE - засоби."
#    (drop the comma right after "CASE")
# -----------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("що є основою CASE, - засоби", $true, $false, $false, $false, $false, $true, 1, $false, "що є основою CASE - засоби", 2) | Out-Null

# -----------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: delete it from its old spot (an empty
#    paragraph before the "Сутності/Відносини" figure) and re-insert it
#    around "Розробкою UML" inside "Розробкою UML керує Object Management
#    Group...".
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$r2 = $d.Content
$found2 = $r2.Find.Execute("Розробкою UML", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $bmRange = $d.Range($r2.Start, $r2.End)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# -----------------------------------------------------------------
# 3) Update the footer's PAGE field cached result from "2" to "7".
# -----------------------------------------------------------------
foreach ($sr in $d.StoryRanges) {
    if ($sr.StoryType -eq 9) {
        foreach ($fld in $sr.Fields) {
            $resStart = $fld.Result.Start
            $resEnd = $fld.Result.End
            $valRange = $sr.Duplicate
            $valRange.Start = $resStart + 2
            $valRange.End = $resEnd + 2
            if ($valRange.Text -eq "2") {
                $valRange.Text = "7"
            }
        }
    }
}
